$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.119
$ws.Range("C14").Value = -12.592
$ws.Range("C21").Value = -12.63
$ws.Range("C23").Value = -12.68
$ws.Range("C25").Value = -11.549
$ws.Range("C26").Value = -12.534
$ws.Range("C29").Value = -12.517
$ws.Range("C53").Value = -12.19
$ws.Range("C57").Value = -13.742
$ws.Range("C59").Value = -12.788
$ws.Range("C69").Value = -10.613
$ws.Range("C79").Value = -12.403
$ws.Range("C83").Value = -13.012
$ws.Range("C91").Value = -12.385
$ws.Range("C93").Value = -10.281
$ws.Range("C103").Value = -12.155
